# Update app to include multi-year boxplots, clean foggy days out of data
#
# The 2001 tower count sheet had a "foggy day" entry (2001-06-15 / "no count")
# that needs to be removed from the dataset so later multi-year boxplot
# analysis isn't skewed by non-count placeholder rows. Deleting the whole
# worksheet row shifts every later observation up by one row, which is
# exactly what the cleaned-up data requires.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 holds the foggy/no-count day (6/15/2001, "no count") - remove it
# entirely so all subsequent rows shift up and the orphaned "no count"
# shared string is dropped from the workbook.
$ws.Rows.Item(5).Delete()

# Leave the selection where the author's session ended up.
$ws.Range("G28").Select()
